# Rename the "wt" / "dcin5" input sheets to make clear they hold
# log2-expression data, then make "dcin5_log2_expression" the active
# (selected/shown) sheet, and scroll the tab strip so
# "wt_log2_expression" is the left-most visible tab.

$wb = $excel.ActiveWorkbook

$wsWt    = $wb.Worksheets.Item("wt")
$wsDcin5 = $wb.Worksheets.Item("dcin5")

$wsWt.Name    = "wt_log2_expression"
$wsDcin5.Name = "dcin5_log2_expression"

# Make dcin5_log2_expression the active/selected tab (was
# "optimization_parameters" before -> activeTab moves 6 -> 3).
$wsDcin5.Activate()
$wsDcin5.Select()

# Scroll the visible tab strip so wt_log2_expression (index 2, 0-based)
# becomes the first shown tab (was dcin5 at index 3 -> wt at index 2).
try {
    $excel.ActiveWindow.ScrollWorkbookTabs(0, 2)
} catch {
}
try {
    $excel.ActiveWindow.ScrollWorkbookTabs(1, 2)
} catch {
}
